{"js": "// Fix a capitalisation typo: \"sende Sie an\" -> \"sende sie an\"\n// (the pronoun \"sie\" was wrongly capitalised).\nconst searchResults = context.document.body.search(\"sende Sie an vorstand@codeforniederrhein.de.\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Target text \"sende Sie an vorstand@codeforniederrhein.de.\" not found.');\n}\n\nsearchResults.items[0].insertText(\n  \"sende sie an vorstand@codeforniederrhein.de.\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# Fix a capitalisation typo: \"sende Sie an\" -> \"sende sie an\"\n# (the pronoun \"sie\" was wrongly capitalised).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"sende Sie an vorstand@codeforniederrhein.de.\"\n$find.Replacement.Text = \"sende sie an vorstand@codeforniederrhein.de.\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Wrap = 1  # wdFindContinue\n$find.Execute([ref]$find.Text, [ref]$find.MatchCase, [ref]$find.MatchWholeWord, $null, $null, $null, $null, $null, $null, [ref]$find.Replacement.Text, 2) | Out-Null\n"}
